$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1 (2)")

# --- Fix the theoretical constant used throughout column Z (rows 2:15) ---
# Old (stale) constant: 0.000000000000107506
# New (corrected) constant: 1.0606413726E-13  (matches the constant already used in column P)
$ws.Range("Z2:Z15").Formula = "=1/(1+(1.0606413726E-13/(9.109E-31*POWER(300000000,2)))*(1-COS(RADIANS(A2))))"

# --- Update the view: scroll/selection moved from L21 to X4 ---
$excel.Goto($ws.Range("K1"), $true)
$ws.Range("X4").Select()
